# Actualización automática 2025-11-28 16:30:09
# Applies the updated sales figures (advisor CASTRO ALCIVAR EDA MARIA)
# across the three report sheets: the per-client/per-group breakdown,
# the monthly sales breakdown, and the monthly-compliance rollup.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" — per-client sales by product group
# ---------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# FERROCENTER CONSTRUCCION ACABADOS FERCONAC SAS gains a 240X80 PORCELANATO sale
$wsGrupo.Range("D27").Value = 1330.56

# JACOME MONCAYO JAVIER ALFONSO gains INODOROS / LAVABOS / PANELES DECORATIVOS
# sales and an increase in PORCELANATO
$wsGrupo.Range("H33").Value = 756.03
$wsGrupo.Range("I33").Value = 73.95
$wsGrupo.Range("K33").Value = 369.98
$wsGrupo.Range("M33").Value = 3406.28

# Row 62 "N de 60" client-count footers move up by one for the touched columns
$wsGrupo.Range("D62").Value = "10 de 60"
$wsGrupo.Range("H62").Value = "4 de 60"
$wsGrupo.Range("I62").Value = "6 de 60"
$wsGrupo.Range("K62").Value = "10 de 60"

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL" — per-client sales by month
# ---------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# noviembre (column F) picks up the same new sales recorded above
$wsMensual.Range("F27").Value = 1330.56
$wsMensual.Range("F33").Value = 4606.24

# Column total for noviembre
$wsMensual.Range("F62").Value = 81636.24000000001

# ---------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" — compliance rollup by product group
# ---------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 23431.83
$wsCumpl.Range("E3").Value = -8606.420000000002
$wsCumpl.Range("F3").Value = 1.580518177912112

# FREGADEROS DE COCINA
$wsCumpl.Range("D4").Value = 1819.6
$wsCumpl.Range("E4").Value = -1030.22
$wsCumpl.Range("F4").Value = 2.305100205224353

# GRIFERIAS
$wsCumpl.Range("D5").Value = 254.89
$wsCumpl.Range("E5").Value = -168.48
$wsCumpl.Range("F5").Value = 2.949774331674575

# INODOROS
$wsCumpl.Range("D6").Value = 5375.05
$wsCumpl.Range("E6").Value = -3770.05
$wsCumpl.Range("F6").Value = 3.348940809968847

# LAVABOS
$wsCumpl.Range("D7").Value = 1467.04
$wsCumpl.Range("E7").Value = -580.328983712426
$wsCumpl.Range("F7").Value = 1.654473636903837

# PANELES DECORATIVOS
$wsCumpl.Range("D10").Value = 8724.57
$wsCumpl.Range("E10").Value = 1191.43
$wsCumpl.Range("F10").Value = 0.8798477208551835

# PORCELANATO
$wsCumpl.Range("D12").Value = 36792.77
$wsCumpl.Range("E12").Value = 13514.23
$wsCumpl.Range("F12").Value = 0.7313648200051682

# PUERTAS DE SEGURIDAD
$wsCumpl.Range("D13").Value = 387.21
$wsCumpl.Range("E13").Value = 723.22665120341
$wsCumpl.Range("F13").Value = 0.3487006661572005

# TOTAL
$wsCumpl.Range("D14").Value = 86213.86
$wsCumpl.Range("E14").Value = 11648.02766749098
$wsCumpl.Range("F14").Value = 0.8809748315190085
